$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string rows appended at the bottom (row 11)
$ws.Range("A11").Value = "Price Option Page check for hints regarding mandatory fields"
$ws.Range("B11").Value = "<HINT Select at least 1 options>"

# Column width changes for columns A and B (drop bestFit, widen)
$ws.Columns.Item(1).ColumnWidth = 59.0
$ws.Columns.Item(2).ColumnWidth = 26.33

# The image anchored to the sheet needs its "to" cell recomputed because
# the widened columns push its right edge into column E (index 4) instead
# of column F (index 5). Re-assert its geometry (in points, 1pt = 12700 EMU)
# using the *exact* original EMU position/size so nothing else shifts.
$shp = $ws.Shapes.Item(1)
$shp.Left = 0
$shp.Top = 197.63511811023622
$shp.Width = 865.8
$shp.Height = 501.2827559055118

$ws.Range("C11").Select()
